# Append the 2025-03-28 price row (row 27) to every price sheet in the
# workbook, carrying forward the prior day's (2025-03-27) value.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-28"

# Sheet name -> new Price value for row 27 (same as the 2025-03-27 value).
$updates = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.21"
    "Cell Topcon 183mm"          = "0.303"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,436"
    "Silver Busbar front-side"   = "8,138"
    "Silver finger front-side"   = "8,188"
    "USD_CNY"                    = "7.2897"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Leading apostrophe forces text so Excel doesn't reinterpret the date
    # string or the numeric-looking price text as a real date/number.
    $ws.Range("A27").Value = "'" + $newDate
    $ws.Range("B27").Value = "'" + $updates[$sheetName]

    # Reset to the default (unstyled) style so no quote-prefix / number
    # formatting is left behind on the new cells, matching the rest of
    # the sheet's plain, style-less cells.
    $ws.Range("A27:B27").Style = "Normal"
}
